# implementacion de filtro para ofertas duplicadas
#
# The "id" column (A) used to hold a random UUID per offer; it now holds
# the numeric LinkedIn job id (the same number embedded in that row's
# "link" URL in column B). Two freshly-seen offers are also appended as
# rows 9-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("offers")

function Set-TextValue($range, [string]$text) {
    # Force the cell to stay a TEXT cell even when the string looks like a
    # pure number (e.g. "4327226302") by entering it with a leading
    # apostrophe (Excel's standard "store as text" quote-prefix), then
    # strip the resulting quote-prefix cell style back to Normal so no
    # stray number formatting is left behind on the cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Column A: replace the UUID ids with the numeric LinkedIn job id already
# present in that row's link (column B).
Set-TextValue $ws.Range("A2") "4327226302"
Set-TextValue $ws.Range("A3") "4327106573"
Set-TextValue $ws.Range("A4") "4343521499"
Set-TextValue $ws.Range("A5") "4343532281"
Set-TextValue $ws.Range("A6") "4338300150"
Set-TextValue $ws.Range("A7") "4338260266"
Set-TextValue $ws.Range("A8") "4291659031"

# Append two new offer rows (9 and 10) found after de-duplicating offers.
Set-TextValue $ws.Range("A9") "4343326779"
$ws.Range("B9").Value = "https://www.linkedin.com/jobs/view/4343326779/"
$ws.Range("C9").Value = "2025-12-13T15:35:28+00:00"
$ws.Range("D9").Value = "«Computer vision»: Remote AI Engineer - HireLATAM y más"
$ws.Range("F9").Value = "En Baufest, nuestra "

Set-TextValue $ws.Range("A10") "4338290130"
$ws.Range("B10").Value = "https://www.linkedin.com/jobs/view/4338290130/"
$ws.Range("C10").Value = "2025-12-13T15:35:28+00:00"
$ws.Range("D10").Value = "«Computer vision»: Remote AI Engineer - HireLATAM y más"
$ws.Range("F10").Value = "Join Our Team`nOowli"
# The embedded newline makes the engine auto-expand the row height; re-fit
# it so row 10 keeps the sheet's default (non-custom) row height, matching
# every other data row.
$ws.Rows.Item(10).AutoFit()
